# Fix some remaining pandas and welly bugs:
# the "survey" sheet's dip/plunge column (D) was off by a constant
# 90-degree offset (the values were stored as the complement of the
# intended dip angle). Correct every data row by adding 90 to the
# existing value, then leave the selection parked on the corrected
# column so the fix is obvious when the workbook is reopened.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

$lastRow = $ws.Cells(1, 1).End(-4121).Row   # xlDown
if ($lastRow -lt 2) { $lastRow = 23 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $old = $cell.Value2
    if ($old -ne $null) {
        $cell.Value2 = $old + 90
    }
}

# Reflect the reviewed range in the sheet's selection.
$ws.Range("D2:D23").Select() | Out-Null
